# Move model from 3.3.1 to 3.4
# - BEPEfCT sheet changes from a single boolean toggle to a per-industry-sector
#   boolean toggle (25 sectors), with a new italic unit header on A1.
# - About sheet gains two new explanatory lines at the bottom.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "About"
$ws2 = $wb.Worksheets.Item(2)   # "BEPEfCT"

# ---------------------------------------------------------------------------
# About sheet: existing lines are unchanged; append two new explanatory rows.
# ---------------------------------------------------------------------------
$ws1.Range("A13").Value = "In the U.S., we exempt agriculture and water and waste process emissions. Generally, "
$ws1.Range("A14").Value = "proposed taxes do not cover these sectors."

# ---------------------------------------------------------------------------
# BEPEfCT sheet: rebuild as a per-sector boolean table.
# ---------------------------------------------------------------------------

# Header row: A1 gets a new italic "Unit" label, B1 keeps its existing text.
$ws2.Range("A1").Value = "Unit: boolean (0 or 1)"
$ws2.Range("A1").Font.Italic = $true
$ws2.Range("B1").Value = "Exempt Process Emissions from Carbon Tax"

# Data rows: one row per industry sector, default 0 except the two exempted
# sectors (agriculture and forestry; water and waste), which are 1.
$sectors = @(
    "agriculture and forestry 01T03",
    "coal mining 05",
    "oil and gas extraction 06",
    "other mining and quarrying 07T08",
    "food beverage and tobacco 10T12",
    "textiles apparel and leather 13T15",
    "wood products 16",
    "pulp paper and printing 17T18",
    "refined petroleum and coke 19",
    "chemicals 20",
    "rubber and plastic products 22",
    "glass and glass products 231",
    "cement and other nonmetallic minerals 239",
    "iron and steel 241",
    "other metals 242",
    "metal products except machinery and vehicles 25",
    "computers and electronics 26",
    "appliances and electrical equipment 27",
    "other machinery 28",
    "road vehicles 29",
    "nonroad vehicles 30",
    "other manufacturing 31T33",
    "energy pipelines and gas processing 352T353",
    "water and waste 36T39",
    "construction 41T43"
)

$exempt = @{
    "agriculture and forestry 01T03" = 1
    "water and waste 36T39" = 1
}

$row = 2
foreach ($sector in $sectors) {
    $value = 0
    if ($exempt.ContainsKey($sector)) {
        $value = $exempt[$sector]
    }
    $ws2.Cells.Item($row, 1).Value = $sector
    $ws2.Cells.Item($row, 2).Value = $value
    $row = $row + 1
}

# Column A is widened to fit the longer sector names.
$ws2.Range("A1").ColumnWidth = 47.18

# Page orientation explicitly set (matches the "About" sheet's existing setup).
$ws2.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Selections / active sheet: "About" stays the active tab (with A15 selected,
# i.e. just below the new content) while "BEPEfCT" remembers B5 selected.
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("B5").Select()

$ws1.Activate()
$ws1.Range("A15").Select()
